# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity, Temperature and mmWave
# sheets, matching the source device export format (all values stored as
# plain text, never auto-converted to dates/times/numbers/percentages).

$wb = $excel.ActiveWorkbook

function Set-LogRow($ws, $row, [string]$date, [string]$timestamp, [string]$hour, [string]$location, [string]$value, [string]$status) {
    $vals = @($date, $timestamp, $hour, $location, $value, $status)
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        # Force text storage so values like dates/times/percentages aren't
        # reinterpreted as numbers, then drop back to the default "Normal"
        # style so no explicit cell format sticks around afterwards.
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$col - 1]
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# PIR sheet: rows 201-213, Bathroom / No Motion / Inactive
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(201, "15:04:17"),
    @(202, "15:04:19"),
    @(203, "15:04:25"),
    @(204, "15:04:29"),
    @(205, "15:04:34"),
    @(206, "15:04:39"),
    @(207, "15:04:44"),
    @(208, "15:04:49"),
    @(209, "15:04:54"),
    @(210, "15:04:59"),
    @(211, "15:05:05"),
    @(212, "15:05:09"),
    @(213, "15:05:14")
)
foreach ($r in $pirRows) {
    $rowNum = $r[0]
    $tstamp = $r[1]
    Set-LogRow $wsPIR $rowNum "2026-01-28" $tstamp "15:00" "Bathroom" "No Motion" "Inactive"
}

# ---------------------------------------------------------------------------
# Humidity sheet: rows 194-203, Bathroom / <percent> / Active
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(194, "15:04:15", "88.1%"),
    @(195, "15:04:24", "88.1%"),
    @(196, "15:04:28", "89.0%"),
    @(197, "15:04:36", "88.9%"),
    @(198, "15:04:48", "88.9%"),
    @(199, "15:04:52", "87.9%"),
    @(200, "15:04:56", "88.8%"),
    @(201, "15:05:04", "87.9%"),
    @(202, "15:05:08", "88.8%"),
    @(203, "15:05:12", "87.9%")
)
foreach ($r in $humidityRows) {
    $rowNum = $r[0]
    $tstamp = $r[1]
    $pct = $r[2]
    Set-LogRow $wsHumidity $rowNum "2026-01-28" $tstamp "15:00" "Bathroom" $pct "Active"
}

# ---------------------------------------------------------------------------
# Temperature sheet: rows 194-203, Bathroom / <celsius> / Active
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(194, "15:04:16", "23.0C"),
    @(195, "15:04:24", "23.0C"),
    @(196, "15:04:28", "23.0C"),
    @(197, "15:04:36", "23.0C"),
    @(198, "15:04:48", "23.0C"),
    @(199, "15:04:52", "23.0C"),
    @(200, "15:04:56", "23.0C"),
    @(201, "15:05:04", "23.0C"),
    @(202, "15:05:08", "23.0C"),
    @(203, "15:05:13", "23.0C")
)
foreach ($r in $temperatureRows) {
    $rowNum = $r[0]
    $tstamp = $r[1]
    $temp = $r[2]
    Set-LogRow $wsTemperature $rowNum "2026-01-28" $tstamp "15:00" "Bathroom" $temp "Active"
}

# ---------------------------------------------------------------------------
# mmWave sheet: rows 9-10, Living Room / Presence value / Status
# ---------------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")
Set-LogRow $wsMmWave 9  "2026-01-28" "15:04:25" "15:00" "Living Room" "No Presence"       "Inactive"
Set-LogRow $wsMmWave 10 "2026-01-28" "15:04:33" "15:00" "Living Room" "Presence Detected" "Active"
